# edit.ps1
# Applies the changes described by the diff:
#  1. Remove the "Meta description" paragraph that originally followed the title.
#  2. Insert a new bold paragraph "Play Book of Relics slot game for free"
#     right before the final (italic "Prompt: ...") paragraph.
#  3. Replace the text of that final italic paragraph with the new
#     meta-description text, keeping its italic formatting.

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph (2nd paragraph) ---
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete() | Out-Null
}

# --- Step 2: insert new bold paragraph just before the last paragraph ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$newParaObj = $d.Paragraphs.Item($count)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Relics slot game for free</w:t></w:r></w:p>'
$newParaObj.Range.InsertXML($newParaXml) | Out-Null

# --- Step 3: replace the text of the final (italic "Prompt: ...") paragraph ---
$oldText = 'Prompt: Create a cartoon-style feature image for "Book of Relics" that features a happy Maya warrior with glasses. The image should be colorful and eye-catching, with the Maya warrior standing in front of an ancient temple or pyramid, holding the Book of Relics in one hand and a handful of gold coins in the other. The background should be a desert landscape, with palm trees and sand dunes visible in the distance. The Maya warrior should be depicted with a big smile on their face, looking excited and confident as they hold their treasures. The image should evoke a sense of adventure, excitement, and the thrill of discovering ancient relics and treasure.'
$newText = 'Read our review of Book of Relics slot game and play for free. Discover its ancient Egypt theme, features and flexible betting range.'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
